# Fixes an error with labeling where we were failing to account for the
# market share affected (20%); drops values from 0.1 to 0.02.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("PPEIdtIL")

# --- Data fix: heating / cooling & ventilation / appliances rows (2, 3, 6) ---
# Column B holds the raw input; Column C is a formula (=B); Column D is a
# second raw input that mirrors B. All three move from 0.1 to 0.02.
$wsData.Range("B2").Value = 0.02
$wsData.Range("D2").Value = 0.02

$wsData.Range("B3").Value = 0.02
$wsData.Range("D3").Value = 0.02

$wsData.Range("B6").Value = 0.02
$wsData.Range("D6").Value = 0.02

# --- View-state: the file was last saved with the About sheet active and a
# block selected; the new save has the PPEIdtIL sheet active with L3 selected.
$wsAbout.Activate()
$wsAbout.Range("E35").Select()

$wsData.Activate()
$wsData.Range("L3").Select()
